$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text change: H3 "Penyelia SRM" -> "Pemimpin Cabang" ---
$ws.Range("H3").Value = "Pemimpin Cabang"

# --- Value changes row 2 ---
$ws.Range("F2").Value = 23320
$ws.Range("M2").Value = 10201872811
$ws.Range("O2").Value = 9448808661
$ws.Range("O2").VerticalAlignment = -4108   # xlCenter -> new style (fontId4 + vertical center)

# --- Value changes row 3 ---
$ws.Range("F3").Style = "Normal"
$ws.Range("F3").Font.Color = 0              # new font (Calibri 11 black, scheme minor)
$ws.Range("F3").Value = 20478

# --- Value changes row 4 ---
$ws.Range("F4").Value = 23320
$ws.Range("M4").Value = 10201872811
$ws.Range("O4").Value = 9448808661
$ws.Range("O4").VerticalAlignment = -4108

# --- Column L width ---
$ws.Columns.Item(12).ColumnWidth = 35.5

# --- Selection ---
$ws.Range("L4").Select()
